$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @(1169, 45534.5, 0.4209, 0.4254, 0.4157, 0.4204, 9198405.800000001),
  @(1170, 45534.58333333334, 0.4205, 0.4225, 0.399, 0.4016, 20799445.4),
  @(1171, 45534.66666666666, 0.4017, 0.4116, 0.4001, 0.4108, 11564426.5),
  @(1172, 45534.75, 0.4107, 0.4213, 0.41, 0.4157, 8969279.4),
  @(1173, 45534.83333333334, 0.4157, 0.4203, 0.4155, 0.4196, 2855969),
  @(1174, 45534.91666666666, 0.4196, 0.4219, 0.4191, 0.421, 3835978.5),
  @(1175, 45535, 0.421, 0.4241, 0.4198, 0.4219, 3330259.8),
  @(1176, 45535.08333333334, 0.4219, 0.4256, 0.4208, 0.4242, 3236084.1),
  @(1177, 45535.16666666666, 0.4242, 0.427, 0.4231, 0.4238, 4195562.6),
  @(1178, 45535.25, 0.4238, 0.4253, 0.4217, 0.4218, 2994097.3),
  @(1179, 45535.33333333334, 0.4218, 0.4236, 0.4181, 0.4194, 3941032.5),
  @(1180, 45535.41666666666, 0.4194, 0.4244, 0.4193, 0.4224, 3791242.3),
  @(1181, 45535.5, 0.4224, 0.4225, 0.4185, 0.4188, 3291590.9),
  @(1182, 45535.58333333334, 0.4189, 0.4215, 0.4181, 0.4214, 4815152.9),
  @(1183, 45535.66666666666, 0.4214, 0.4215, 0.4173, 0.4194, 6736648.2),
  @(1184, 45535.75, 0.4194, 0.4233, 0.4179, 0.4191, 2274342.4),
  @(1185, 45535.83333333334, 0.4191, 0.4198, 0.4179, 0.4194, 1103924.6),
  @(1186, 45535.91666666666, 0.4209, 0.4209, 0.4188, 0.4196, 1196986.3),
  @(1187, 45536, 0.4196, 0.4203, 0.4146, 0.4148, 1598104.6),
  @(1188, 45536.08333333334, 0.4148, 0.4164, 0.4106, 0.4148, 3975380),
  @(1189, 45536.16666666666, 0.4148, 0.4149, 0.4057, 0.4064, 4881467.4),
  @(1190, 45536.25, 0.4065, 0.4107, 0.4059, 0.4106, 5683429.1),
  @(1191, 45536.33333333334, 0.4107, 0.4122, 0.4091, 0.4105, 3719244.8),
  @(1192, 45536.41666666666, 0.4104, 0.4129, 0.4077, 0.4127, 4360750.7),
  @(1193, 45536.5, 0.4126, 0.4135, 0.405, 0.4077, 4037056.5),
  @(1194, 45536.58333333334, 0.4076, 0.412, 0.4014, 0.4108, 8377985.1),
  @(1195, 45536.66666666666, 0.4108, 0.4124, 0.4044, 0.4058, 4665393.9),
  @(1196, 45536.75, 0.4057, 0.4136, 0.4055, 0.4133, 5024926.4),
  @(1197, 45536.83333333334, 0.4132, 0.4162, 0.408, 0.4105, 3731954.8),
  @(1198, 45536.91666666666, 0.4105, 0.4108, 0.3967, 0.4023, 7960578.4),
  @(1199, 45537, 0.4022, 0.4071, 0.4012, 0.4046, 3706815.5),
  @(1200, 45537.08333333334, 0.4046, 0.409, 0.4042, 0.4072, 3968256.8),
  @(1201, 45537.16666666666, 0.4073, 0.4083, 0.4005, 0.4067, 6684150.4),
  @(1202, 45537.25, 0.4067, 0.4083, 0.4002, 0.4015, 5748045.6),
  @(1203, 45537.33333333334, 0.4015, 0.406, 0.3987, 0.4039, 15704896.6),
  @(1204, 45537.41666666666, 0.4038, 0.4093, 0.4022, 0.4055, 9733844.300000001),
  @(1205, 45537.5, 0.4056, 0.4061, 0.3982, 0.4, 9536347.5),
  @(1206, 45537.58333333334, 0.4001, 0.4045, 0.3956, 0.4034, 12359047.8),
  @(1207, 45537.66666666666, 0.4033, 0.4072, 0.4016, 0.4069, 4138916.1),
  @(1208, 45537.75, 0.4069, 0.4091, 0.4051, 0.4066, 2906117.5),
  @(1209, 45537.83333333334, 0.4067, 0.415, 0.4058, 0.4116, 5963414.7),
  @(1210, 45537.91666666666, 0.4116, 0.4129, 0.4104, 0.4108, 1995340.8),
  @(1211, 45538, 0.4108, 0.4108, 0.4075, 0.4099, 3047355.4),
  @(1212, 45538.08333333334, 0.4099, 0.415, 0.4091, 0.4104, 8412503),
  @(1213, 45538.16666666666, 0.4103, 0.4128, 0.4094, 0.4128, 5064586.9),
  @(1214, 45538.25, 0.4127, 0.4129, 0.4099, 0.4126, 6430799.6),
  @(1215, 45538.33333333334, 0.4127, 0.4128, 0.4057, 0.4073, 5810397.4),
  @(1216, 45538.41666666666, 0.4073, 0.4107, 0.4073, 0.4079, 5582405.6),
  @(1217, 45538.5, 0.408, 0.4095, 0.3998, 0.4027, 8381071.5),
  @(1218, 45538.58333333334, 0.4027, 0.4127, 0.4002, 0.4058, 15357562.4),
  @(1219, 45538.66666666666, 0.4058, 0.4072, 0.3996, 0.4002, 7915846.4),
  @(1220, 45538.75, 0.4003, 0.405, 0.3982, 0.4031, 2902158.4),
  @(1221, 45538.83333333334, 0.4032, 0.4059, 0.4015, 0.4027, 7009858.3),
  @(1222, 45538.91666666666, 0.4027, 0.4037, 0.3965, 0.3978, 4699639.8),
  @(1223, 45539, 0.3978, 0.4003, 0.3781, 0.3906, 18717904.5),
  @(1224, 45539.08333333334, 0.3907, 0.3923, 0.3886, 0.3893, 7502502.5),
  @(1225, 45539.16666666666, 0.3893, 0.3906, 0.3839, 0.3845, 6627530.4),
  @(1226, 45539.25, 0.3844, 0.393, 0.3843, 0.3896, 8173478.1),
  @(1227, 45539.33333333334, 0.3895, 0.3895, 0.35, 0.3744, 52295760.2),
  @(1228, 45539.41666666666, 0.3743, 0.3745, 0.3672, 0.3738, 11442791),
  @(1229, 45539.5, 0.3738, 0.3769, 0.3714, 0.3748, 6426266.6),
  @(1230, 45539.58333333334, 0.3747, 0.385, 0.3733, 0.3845, 7979953.1),
  @(1231, 45539.66666666666, 0.3847, 0.3921, 0.3806, 0.3834, 11030555.1),
  @(1232, 45539.75, 0.3833, 0.3838, 0.3799, 0.3812, 1639687.3)
)

foreach ($r in $rows) {
  $rowNum = $r[0]
  $ws.Cells.Item($rowNum, 1).Value2 = $r[1]
  $ws.Cells.Item($rowNum, 2).Value2 = $r[2]
  $ws.Cells.Item($rowNum, 3).Value2 = $r[3]
  $ws.Cells.Item($rowNum, 4).Value2 = $r[4]
  $ws.Cells.Item($rowNum, 5).Value2 = $r[5]
  $ws.Cells.Item($rowNum, 6).Value2 = $r[6]
}

# Apply the same formatting (number format / font / border / alignment)
# used in column A for existing rows, to the A column of all new rows.
$srcFormat = $ws.Range("A1168")
$dstFormat = $ws.Range("A1169:A1232")
$srcFormat.Copy()
$dstFormat.PasteSpecial(-4122)
